$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gradient")

# Update the BottomLeft and BottomRight gradient input values
$ws.Range("B4").Value = 20
$ws.Range("B5").Value = 50

# Update the active selection on the sheet to B4
$ws.Range("B4").Select()
